# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the existing "Late" / "Outstanding" columns one place to
# the right, then make "Repayment Schedule" the active sheet with the
# newly-revealed column selected.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsTransactions = $wb.Worksheets.Item("Transactions")

# Insert a blank column before column N (pushes N->O, O->P, P->Q, ...).
# Excel's default "Insert" behaviour copies the format (incl. width) of the
# column to the left, so match column M's width on the newly created column.
$wsSchedule.Columns("N:N").Insert()
$wsSchedule.Columns("N:N").ColumnWidth = 10.5

# Activate the Repayment Schedule sheet and select the new data cell.
$wsSchedule.Activate()
$wsSchedule.Range("Q8").Select()

# The Transactions sheet is no longer the active tab.
$wsTransactions.Range("F8").Select()

$wsSchedule.Activate()
